$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("sigma_010")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 27.62255115134589
$ws.Cells.Item(2, 3).Value = 29.93256963034784
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 27.60707422963498
$ws.Cells.Item(3, 3).Value = 29.94122657403026
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 27.62306503157641
$ws.Cells.Item(4, 3).Value = 29.97416284828283
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 27.60616732129732
$ws.Cells.Item(5, 3).Value = 29.95164715338845
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 27.60213676552733
$ws.Cells.Item(6, 3).Value = 29.97624006148319
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 27.63458644029293
$ws.Cells.Item(7, 3).Value = 29.97546539490484
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 27.59956342976039
$ws.Cells.Item(8, 3).Value = 29.96223137930503
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 27.6303883290899
$ws.Cells.Item(9, 3).Value = 29.97985750529905
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 27.58007418012268
$ws.Cells.Item(10, 3).Value = 29.96640279206242
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 27.64660357968289
$ws.Cells.Item(11, 3).Value = 29.97885263396768
$ws.Cells.Item(12, 2).Value = 27.61522104583307
$ws.Cells.Item(12, 3).Value = 29.96386559730716

$ws = $wb.Worksheets.Item("sigma_025")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 19.73626406716477
$ws.Cells.Item(2, 3).Value = 26.77869325068998
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 19.71376546007914
$ws.Cells.Item(3, 3).Value = 26.7432427556712
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 19.73180231998506
$ws.Cells.Item(4, 3).Value = 26.79282029396019
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 19.68864522651385
$ws.Cells.Item(5, 3).Value = 26.76271522637571
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 19.71758351508478
$ws.Cells.Item(6, 3).Value = 26.80726110502393
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 19.69758103777424
$ws.Cells.Item(7, 3).Value = 26.71545803224972
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 19.7008796529013
$ws.Cells.Item(8, 3).Value = 26.7691591662143
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 19.70283581960065
$ws.Cells.Item(9, 3).Value = 26.81484375085085
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 19.7153506279143
$ws.Cells.Item(10, 3).Value = 26.77506201169315
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 19.70090055591832
$ws.Cells.Item(11, 3).Value = 26.79691322174576
$ws.Cells.Item(12, 2).Value = 19.71056082829364
$ws.Cells.Item(12, 3).Value = 26.77561688144748

$ws = $wb.Worksheets.Item("sigma_050")
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 14.70839733236271
$ws.Cells.Item(2, 3).Value = 21.40486654949476
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 14.71576074376576
$ws.Cells.Item(3, 3).Value = 21.36575608612224
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 14.71438766487775
$ws.Cells.Item(4, 3).Value = 21.40588918474436
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 14.72559432362701
$ws.Cells.Item(5, 3).Value = 21.43752019339114
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 14.72253718223117
$ws.Cells.Item(6, 3).Value = 21.39503518784375
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 14.69908148781415
$ws.Cells.Item(7, 3).Value = 21.38167857046971
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 14.70370540700552
$ws.Cells.Item(8, 3).Value = 21.33820974984896
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 14.70645596129224
$ws.Cells.Item(9, 3).Value = 21.38578734474979
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 14.73039107767457
$ws.Cells.Item(10, 3).Value = 21.3767082833361
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 14.71580268507553
$ws.Cells.Item(11, 3).Value = 21.33780149540972
$ws.Cells.Item(12, 2).Value = 14.71421138657264
$ws.Cells.Item(12, 3).Value = 21.38292526454105
